# Commit: "fix typos in sample arc, update playground.fsx"
#
# The workbook has two sheets: "isa_study" (ISA metadata) and "Sheet1"
# (the growth-conditions data sheet, which gets renamed to "Growth").
# On the "Growth" sheet, column DQ ("Sample Name") holds three sample
# names that were typo'd: "C1" and "C2" (with "C2" duplicated across two
# rows). Fix them to the unique, correctly-named "CC1" / "CC2" / "CC3",
# rename the sheet, and leave the selection on the last corrected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename "Sheet1" -> "Growth"
$ws.Name = "Growth"

# Fix the sample-name typos in column DQ (rows 2-4)
$ws.Range("DQ2").Value = "CC1"
$ws.Range("DQ3").Value = "CC2"
$ws.Range("DQ4").Value = "CC3"

# Leave the sheet active with the last corrected cell selected
$ws.Activate() | Out-Null
$ws.Range("DQ4").Select() | Out-Null
